$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftmost (index/style) column; remaining columns shift left.
$ws.Columns("A").Delete()

# Fix header text typo: MODEL_CONDITION -> MODELCONDITION (now in column D).
$ws.Range("D1").Value = "MODELCONDITION"
